# Weekly update: insert the latest week's price-report row for
# "Poroto verde" (Terminal La Palmera de La Serena) at the top of the
# data block (row 71), pushing every existing row down by one.
#
# This mirrors how the source report is produced each week: a brand new
# observation is prepended and the previously-existing rows simply slide
# down (row 71 -> 72, 72 -> 73, ... 127 -> 128), with the sheet's used
# range growing from A1:R127 to A1:R128.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 71..127 down to 72..128, leaving a blank row 71 behind.
$ws.Rows.Item(71).Insert()

# Fill the newly-opened row 71 with this week's reading.
$ws.Range("A71").Value = 8
$ws.Range("B71").Value = "Terminal La Palmera de La Serena"
$ws.Range("C71").Value = "Coquimbo"
$ws.Range("D71").Value = 44484
$ws.Range("E71").Value = 4
$ws.Range("F71").Value = 100112031
$ws.Range("G71").Value = "Poroto verde"
$ws.Range("H71").Value = "Magnum"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 520
$ws.Range("K71").Value = 35000
$ws.Range("L71").Value = 36000
$ws.Range("M71").Value = 35500
$ws.Range("N71").Value = "$/malla 25 kilos"
$ws.Range("O71").Value = "Perú"
$ws.Range("P71").Value = 1420
$ws.Range("Q71").Value = 25
$ws.Range("R71").Value = "Hortaliza"
